# Weekly update: a new price-report row is inserted at the top of the
# data (row 2), pushing all existing data rows down by one. The sheet's
# dimension grows from A1:R20 to A1:R21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (shifts rows 2..20 -> 3..21)
$ws.Rows.Item(2).Insert()

# The inserted row inherits the header row's (bold/centered) formatting from
# Excel's default "format from above" behavior; reset it back to the plain
# style used by the rest of the data rows.
$ws.Range("A2:R2").ClearFormats()

# Column D carries a date number format on every data row.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row with this week's report.
$ws.Range("D2").Value = Get-Date -Year 2022 -Month 5 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112026
$ws.Range("G2").Value = "Haba"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = 1250
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 1250
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
